$d = $word.ActiveDocument

# --- Change 1: merge the "...are 5 / star / versus..." runs into one run
#     (removes the grammar-check proofErr markers around "star"), while
#     keeping the preceding "... whether " run intact/separate.
$rng = $d.Content
$rng.Find.Execute(
    "there is any bias in the reviews given by paid reviewers (vine). To do this, we’ll look at the percentage of vine reviews that are 5 star versus the number of non-vine reviews that give the same rating.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "there is any bias in the reviews given by paid reviewers (vine). To do this, we’ll look at the percentage of vine reviews that are 5 star versus the number of non-vine reviews that give the same rating.",
    2
) | Out-Null

# The replace above also merges into the preceding run (" ... whether ").
# Re-establish the run boundary right before "there is any bias" so that
# run stays separate, matching the target structure.
$rng3 = $d.Content
$rng3.Find.Execute(
    "there is any bias in the reviews given by paid reviewers (vine). To do this, we’ll look at the percentage of vine reviews that are 5 star versus the number of non-vine reviews that give the same rating.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
) | Out-Null
$rng3.Bold = 1
$rng3.Bold = 0

# --- Change 2: merge "The results of this analysis were / actually compelling / ... / was / ..." runs into one run
#     (removes the grammar-check proofErr markers around "actually compelling" and "was")
$rng = $d.Content
$rng.Find.Execute(
    "The results of this analysis were actually compelling in a very unexpected way. Upon filtering the data to remove any entries where the total votes was less than 20, the results returned absolutely no vine reviews within the dataset:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The results of this analysis were actually compelling in a very unexpected way. Upon filtering the data to remove any entries where the total votes was less than 20, the results returned absolutely no vine reviews within the dataset:",
    2
) | Out-Null

# --- Change 3: merge "...before they are / released / and similar tactics..." runs into one run
#     (removes the grammar-check proofErr markers around "released")
$rng = $d.Content
$rng.Find.Execute(
    "In the modern age of the internet and product reviews, it is easy for companies and individuals to review just about anything, even if they don’t have experience with the thing they are reviewing. Bots or agencies creating fake accounts to boost the perception of a product, users protesting products or movies by review-bombing things before they are released and similar tactics can make the reliability of reviews less than ideal. To combat this, some companies such as amazon have included “verified purchase” tags on reviews where they can confirm that a product was purchased by the reviewer. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In the modern age of the internet and product reviews, it is easy for companies and individuals to review just about anything, even if they don’t have experience with the thing they are reviewing. Bots or agencies creating fake accounts to boost the perception of a product, users protesting products or movies by review-bombing things before they are released and similar tactics can make the reliability of reviews less than ideal. To combat this, some companies such as amazon have included “verified purchase” tags on reviews where they can confirm that a product was purchased by the reviewer. ",
    2
) | Out-Null

# --- Change 4: append a new sentence (new run) at the end of the "likely reliable." paragraph
$rng = $d.Content
$rng.Find.Execute("are therefore likely reliable.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertAfter(" A secondary analysis should be done on one star reviews to determine whether the consistency holds true at the other end of the rating spectrum.")
